# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 53 (pushing the existing
# rows 53-95 down to 54-96), then populate the new row with the
# new "Coco" (coconut) market observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 53..95 down by one row, making room for the new record.
$ws.Rows.Item(53).Insert()

# Fill in the newly inserted row 53 with the new observation.
$ws.Cells.Item(53, 1).Value2 = 6
$ws.Cells.Item(53, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(53, 3).Value2 = "Metropolitana"
$ws.Cells.Item(53, 4).Value2 = 45240
$ws.Cells.Item(53, 5).Value2 = 13
$ws.Cells.Item(53, 6).Value2 = "Fruta"
$ws.Cells.Item(53, 7).Value2 = 100108
$ws.Cells.Item(53, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(53, 9).Value2 = 100108007
$ws.Cells.Item(53, 10).Value2 = "Coco"
$ws.Cells.Item(53, 11).Value2 = "Sin especificar"
$ws.Cells.Item(53, 12).Value2 = "Primera"
$ws.Cells.Item(53, 13).Value2 = 50
$ws.Cells.Item(53, 14).Value2 = 30000
$ws.Cells.Item(53, 15).Value2 = 30000
$ws.Cells.Item(53, 16).Value2 = 30000
$ws.Cells.Item(53, 17).Value2 = "$/malla 20 unidades"
$ws.Cells.Item(53, 18).Value2 = "Ecuador"
$ws.Cells.Item(53, 19).Value2 = 1500
$ws.Cells.Item(53, 20).Value2 = 20
